# update new orleans xlsx files
#
# 1) hotel_info gains a new "State" column, inserted right after
#    "Hotel_Name" (i.e. before "City"), populated with "Louisiana" for
#    the existing data row.
# 2) The two sheet tabs are reordered so "review_info" comes before
#    "hotel_info".

$wb = $excel.ActiveWorkbook

$hotelInfo = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")

# --- 1. Insert the new "State" column into hotel_info ---------------------
# Hotel_Name is column B, City is column C -- insert a fresh column C so
# everything from the old "City" onward shifts one column to the right,
# leaving Hotel_Name (column B) untouched.
$hotelInfo.Range("C1").EntireColumn.Insert()

$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"

# --- 2. Reorder the sheet tabs: review_info, then hotel_info --------------
$hotelInfo.Move($null, $reviewInfo)

Write-Output "done"
